$wb = $excel.ActiveWorkbook

# Sheet "OFF" - update row 3 (label "R") with Week 16 simulated values
$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("B3").Value = 480
$wsOFF.Range("C3").Value = 343
$wsOFF.Range("D3").Value = 120
$wsOFF.Range("E3").Value = 65
$wsOFF.Range("F3").Value = 8

# Sheet "DEF" - update row 3 (label "R") with Week 16 simulated values
$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("B3").Value = 412
$wsDEF.Range("C3").Value = 287
$wsDEF.Range("D3").Value = 106
$wsDEF.Range("E3").Value = 57
$wsDEF.Range("F3").Value = 6
